# Finishing to config the Logging Process
# - Set the "MaxRetryNumber" value (B2) on the Constants sheet to 3.
# - Update the view so the selected/active cell is B5 (scrolled so row 1
#   stays visible with column B at the left edge).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")
$ws.Activate()

# Update the MaxRetryNumber value from 0 to 3
$ws.Range("B2").Value = 3

# Move the selection/active cell to B5 and adjust the visible top-left cell
$ws.Range("B5").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2
